$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values scraped on the latest GitHub Actions run.
# Numeric-looking text values (e.g. "0.0300", "81.50") must be forced to
# remain text so trailing/leading zeros and exact formatting survive -
# otherwise Excel auto-converts them to real numbers.

$ws.Range("D2").Value = "42.739.55"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "2.546.14"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.79%  "
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").Value = "2.938.66"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("E15").Value = "  +5.41%  "
$ws.Range("D16").Value = "2.514.85"
$ws.Range("E16").Value = "  -4.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.840"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "42.766.67"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").Value = "0.0₃0957"
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "246.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.69%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -2.29%  "
$ws.Range("E30").Value = "  -2.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0806"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("E34").Value = "  -2.99%  "
$ws.Range("E35").Value = "  -3.02%  "
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.53%  "
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.86%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.78%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0300"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").Value = "1.984.89"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("D48").Value = "2.792.91"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "81.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("E51").Value = "  -1.64%  "
